$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.885.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.515.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.63%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.514.28'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.64%  '

$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.82'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.383'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.48%  '

$ws.Range("E13").Value = '  +3.23%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.536.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.56%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.90%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000180'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.56%  '

$ws.Range("E17").Value = '  +1.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.028.11'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.84%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.568'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.24%  '

$ws.Range("E24").Value = '  +3.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.06%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000111'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.40%  '

$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("E30").Value = '  +2.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.90%  '

$ws.Range("E32").Value = '  +3.54%  '

$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +16.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.142'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '168.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0795'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.816'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +16.33%  '

$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.385.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '303.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +13.54%  '
